# Updates crypto price/volume snapshot cells (and restores the correct
# Coin/Link pairing for rows 17-18) to match the refreshed data pull.
#
# Price values that look like plain decimals (e.g. "241.66", "1.00") are
# written with a leading apostrophe so Excel keeps them as literal text
# (matching the original inline-string cells) instead of silently
# converting them to numbers and dropping significant trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.079.72'
$ws.Range('E2').Value = '  -1.27%  '

$ws.Range('D3').Value = '1.986.51'
$ws.Range('E3').Value = '  -2.49%  '

$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').Value = "'241.66"
$ws.Range('E5').Value = '  -6.19%  '

$ws.Range('D6').Value = "'0.603"
$ws.Range('E6').Value = '  -3.49%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = "'55.11"
$ws.Range('E8').Value = '  -4.67%  '

$ws.Range('D9').Value = "'0.371"
$ws.Range('E9').Value = '  -4.17%  '

$ws.Range('D10').Value = "'58.29"
$ws.Range('E10').Value = '  +1.67%  '

$ws.Range('D11').Value = "'0.0754"
$ws.Range('E11').Value = '  -5.42%  '

$ws.Range('E12').Value = '  -4.68%  '

$ws.Range('D13').Value = '2.278.29'
$ws.Range('E13').Value = '  -2.49%  '

$ws.Range('D14').Value = "'14.05"
$ws.Range('E14').Value = '  -5.47%  '

$ws.Range('D15').Value = "'20.83"
$ws.Range('E15').Value = '  -3.00%  '

$ws.Range('D16').Value = "'0.756"
$ws.Range('E16').Value = '  -8.09%  '

$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = "'5.04"
$ws.Range('E17').Value = '  -6.18%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '1.981.90'
$ws.Range('E18').Value = '  -2.81%  '

$ws.Range('D19').Value = '36.942.06'
$ws.Range('E19').Value = '  -1.47%  '

$ws.Range('D20').Value = "'68.06"
$ws.Range('E20').Value = '  -2.98%  '

$ws.Range('D21').Value = '0.0₃0809'
$ws.Range('E21').Value = '  -5.58%  '

$ws.Range('D22').Value = "'227.71"
$ws.Range('E22').Value = '  -0.80%  '

$ws.Range('D23').Value = "'4.97"
$ws.Range('E23').Value = '  -4.93%  '

$ws.Range('E24').Value = '  +0.14%  '

$ws.Range('D25').Value = "'2.41"
$ws.Range('E25').Value = '  -9.17%  '

$ws.Range('D26').Value = "'2.35"
$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').Value = "'160.81"
$ws.Range('E27').Value = '  -1.79%  '

$ws.Range('D28').Value = "'8.65"
$ws.Range('E28').Value = '  -5.55%  '

$ws.Range('D29').Value = "'19.12"
$ws.Range('E29').Value = '  -4.76%  '

$ws.Range('D30').Value = "'0.122"
$ws.Range('E30').Value = '  -11.51%  '

$ws.Range('D31').Value = "'1.31"
$ws.Range('E31').Value = '  -3.55%  '

$ws.Range('E32').Value = '  -3.45%  '

$ws.Range('D33').Value = "'4.40"
$ws.Range('E33').Value = '  -7.15%  '

$ws.Range('D34').Value = "'0.0610"
$ws.Range('E34').Value = '  -8.74%  '

$ws.Range('D35').Value = "'4.20"
$ws.Range('E35').Value = '  -7.33%  '

$ws.Range('D36').Value = "'2.32"
$ws.Range('E36').Value = '  -6.14%  '

$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  -0.17%  '

$ws.Range('D38').Value = "'1.79"
$ws.Range('E38').Value = '  -1.50%  '

$ws.Range('D39').Value = "'3.30"
$ws.Range('E39').Value = '  -4.86%  '

$ws.Range('D40').Value = "'5.22"
$ws.Range('E40').Value = '  -3.15%  '

$ws.Range('D41').Value = "'3.10"
$ws.Range('E41').Value = '  +2.13%  '

$ws.Range('D42').Value = '1.430.56'
$ws.Range('E42').Value = '  +2.21%  '

$ws.Range('D43').Value = "'0.0203"
$ws.Range('E43').Value = '  -6.22%  '

$ws.Range('D44').Value = "'1.12"
$ws.Range('E44').Value = '  -6.32%  '

$ws.Range('D45').Value = "'0.0877"
$ws.Range('E45').Value = '  -9.45%  '

$ws.Range('D46').Value = "'88.14"
$ws.Range('E46').Value = '  -3.49%  '

$ws.Range('D47').Value = "'15.19"
$ws.Range('E47').Value = '  -6.52%  '

$ws.Range('D48').Value = "'1.00"
$ws.Range('E48').Value = '  -5.08%  '

$ws.Range('D49').Value = "'2.88"
$ws.Range('E49').Value = '  +0.43%  '

$ws.Range('D50').Value = "'3.66"
$ws.Range('E50').Value = '  +13.62%  '

$ws.Range('D51').Value = "'6.65"
$ws.Range('E51').Value = '  -10.19%  '
